$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the dct:modified timestamp in B20
$ws.Range("B20").Value = "2022-06-04T21:55:11+00:00"

# Update subject labels in column B, rows 23-57 (subject 1 .. subject 35)
for ($row = 23; $row -le 57; $row++) {
    $n = $row - 22
    $ws.Range("B$row").Value = "subject $n"
}
